# Update "想去人数" (interested-count) figures for two exhibition events.
# - Sheet "展览" (Exhibition): row 3 and row 4, column F
# - Sheet "全部类型" (All Types), which mirrors the same two events: row 5 and row 6, column F

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1228
$wsExhibition.Range("F4").Value = 2714

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F5").Value = 1228
$wsAllTypes.Range("F6").Value = 2714
